$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '23.509.32'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('D2').Style = 'Normal'

$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('E2').NumberFormat = 'General'
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.648.75'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E3').NumberFormat = 'General'
$ws.Range('E3').Style = 'Normal'

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9991'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('D4').Style = 'Normal'

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.33%  '
$ws.Range('E4').NumberFormat = 'General'
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9999'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.15%  '
$ws.Range('E5').NumberFormat = 'General'
$ws.Range('E5').Style = 'Normal'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '300.16'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.79%  '
$ws.Range('E6').NumberFormat = 'General'
$ws.Range('E6').Style = 'Normal'

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3788'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('D7').Style = 'Normal'

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.22%  '
$ws.Range('E7').NumberFormat = 'General'
$ws.Range('E7').Style = 'Normal'

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '50.46'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('D8').Style = 'Normal'

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.30%  '
$ws.Range('E8').NumberFormat = 'General'
$ws.Range('E8').Style = 'Normal'

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3502'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = 'Normal'

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.62%  '
$ws.Range('E9').NumberFormat = 'General'
$ws.Range('E9').Style = 'Normal'

$ws.Range('B10').NumberFormat = '@'
$ws.Range('B10').Value = 'Polygon'
$ws.Range('B10').NumberFormat = 'General'
$ws.Range('B10').Style = 'Normal'

$ws.Range('C10').NumberFormat = '@'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('C10').NumberFormat = 'General'
$ws.Range('C10').Style = 'Normal'

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.223'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.22%  '
$ws.Range('E10').NumberFormat = 'General'
$ws.Range('E10').Style = 'Normal'

$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('B11').NumberFormat = 'General'
$ws.Range('B11').Style = 'Normal'

$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('C11').NumberFormat = 'General'
$ws.Range('C11').Style = 'Normal'

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08060'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.61%  '
$ws.Range('E11').NumberFormat = 'General'
$ws.Range('E11').Style = 'Normal'

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9991'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('D12').Style = 'Normal'

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.35%  '
$ws.Range('E12').NumberFormat = 'General'
$ws.Range('E12').Style = 'Normal'

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.13'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.09%  '
$ws.Range('E13').NumberFormat = 'General'
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.309'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.55%  '
$ws.Range('E14').NumberFormat = 'General'
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.260'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('D15').Style = 'Normal'

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.11%  '
$ws.Range('E15').NumberFormat = 'General'
$ws.Range('E15').Style = 'Normal'

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001212'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.70%  '
$ws.Range('E16').NumberFormat = 'General'
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.647.79'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('D17').Style = 'Normal'

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.53%  '
$ws.Range('E17').NumberFormat = 'General'
$ws.Range('E17').Style = 'Normal'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '95.35'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.27%  '
$ws.Range('E18').NumberFormat = 'General'
$ws.Range('E18').Style = 'Normal'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06993'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.20%  '
$ws.Range('E19').NumberFormat = 'General'
$ws.Range('E19').Style = 'Normal'

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.625'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('D20').Style = 'Normal'

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.73%  '
$ws.Range('E20').NumberFormat = 'General'
$ws.Range('E20').Style = 'Normal'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.45'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.01%  '
$ws.Range('E21').NumberFormat = 'General'
$ws.Range('E21').Style = 'Normal'

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('E22').NumberFormat = 'General'
$ws.Range('E22').Style = 'Normal'

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.45'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.67%  '
$ws.Range('E23').NumberFormat = 'General'
$ws.Range('E23').Style = 'Normal'

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '23.502.10'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('E24').NumberFormat = 'General'
$ws.Range('E24').Style = 'Normal'

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.418'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.79%  '
$ws.Range('E25').NumberFormat = 'General'
$ws.Range('E25').Style = 'Normal'

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.023'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.28%  '
$ws.Range('E26').NumberFormat = 'General'
$ws.Range('E26').Style = 'Normal'

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.08'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.64%  '
$ws.Range('E27').NumberFormat = 'General'
$ws.Range('E27').Style = 'Normal'

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '152.00'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = 'Normal'

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.66%  '
$ws.Range('E28').NumberFormat = 'General'
$ws.Range('E28').Style = 'Normal'

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.181'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = 'Normal'

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.10%  '
$ws.Range('E29').NumberFormat = 'General'
$ws.Range('E29').Style = 'Normal'

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '131.65'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.70%  '
$ws.Range('E30').NumberFormat = 'General'
$ws.Range('E30').Style = 'Normal'

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.827.35'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.64%  '
$ws.Range('E31').NumberFormat = 'General'
$ws.Range('E31').Style = 'Normal'

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.893'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('D32').Style = 'Normal'

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.34%  '
$ws.Range('E32').NumberFormat = 'General'
$ws.Range('E32').Style = 'Normal'

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.137'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = 'Normal'

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -5.06%  '
$ws.Range('E33').NumberFormat = 'General'
$ws.Range('E33').Style = 'Normal'

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -7.67%  '
$ws.Range('E34').NumberFormat = 'General'
$ws.Range('E34').Style = 'Normal'

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9906'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = 'Normal'

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -6.45%  '
$ws.Range('E35').NumberFormat = 'General'
$ws.Range('E35').Style = 'Normal'

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02703'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = 'Normal'

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.61%  '
$ws.Range('E36').NumberFormat = 'General'
$ws.Range('E36').Style = 'Normal'

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.08789'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('E37').NumberFormat = 'General'
$ws.Range('E37').Style = 'Normal'

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.934'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('D38').Style = 'Normal'

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -3.05%  '
$ws.Range('E38').NumberFormat = 'General'
$ws.Range('E38').Style = 'Normal'

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2424'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = 'Normal'

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.96%  '
$ws.Range('E39').NumberFormat = 'General'
$ws.Range('E39').Style = 'Normal'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06809'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.93%  '
$ws.Range('E40').NumberFormat = 'General'
$ws.Range('E40').Style = 'Normal'

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '12.90'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.79%  '
$ws.Range('E41').NumberFormat = 'General'
$ws.Range('E41').Style = 'Normal'

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6899'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.39%  '
$ws.Range('E42').NumberFormat = 'General'
$ws.Range('E42').Style = 'Normal'

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.294'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('D43').Style = 'Normal'

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -3.14%  '
$ws.Range('E43').NumberFormat = 'General'
$ws.Range('E43').Style = 'Normal'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.62'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = 'Normal'

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.27%  '
$ws.Range('E44').NumberFormat = 'General'
$ws.Range('E44').Style = 'Normal'

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9995'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('D45').Style = 'Normal'

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('E45').NumberFormat = 'General'
$ws.Range('E45').Style = 'Normal'

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6400'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('D46').Style = 'Normal'

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.82%  '
$ws.Range('E46').NumberFormat = 'General'
$ws.Range('E46').Style = 'Normal'

$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('B47').NumberFormat = 'General'
$ws.Range('B47').Style = 'Normal'

$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('C47').NumberFormat = 'General'
$ws.Range('C47').Style = 'Normal'

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.247'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = 'Normal'

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.46%  '
$ws.Range('E47').NumberFormat = 'General'
$ws.Range('E47').Style = 'Normal'

$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('B48').NumberFormat = 'General'
$ws.Range('B48').Style = 'Normal'

$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('C48').NumberFormat = 'General'
$ws.Range('C48').Style = 'Normal'

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.926'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.72%  '
$ws.Range('E48').NumberFormat = 'General'
$ws.Range('E48').Style = 'Normal'

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '127.10'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('D49').Style = 'Normal'

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.88%  '
$ws.Range('E49').NumberFormat = 'General'
$ws.Range('E49').Style = 'Normal'

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07674'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.84%  '
$ws.Range('E50').NumberFormat = 'General'
$ws.Range('E50').Style = 'Normal'

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.238'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.98%  '
$ws.Range('E51').NumberFormat = 'General'
$ws.Range('E51').Style = 'Normal'
